$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing rows whose scheme label shifted position because three new
# "Spiral-*" schemes (plus the existing "Gaussian-Quadrature" entry) are now
# inserted between "Ring Perpendicular to TD" and "NoRotation-tilt60deg".
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("B16").Value = "Rotation-60detTilt"

# Three brand-new rows (17-19) holding the schemes pushed off the end of the
# original table: HexGrid-90degTilt5degRes / 22p5degRes / 60degTilt5degRes.
$schemes17to19 = @("HexGrid-90degTilt5degRes", "HexGrid-90degTilt22p5degRes", "HexGrid-60degTilt5degRes")
$ws.Range("A16").Copy()
for ($i = 0; $i -lt 3; $i++) {
    $r = 17 + $i
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = 15 + $i
    $ws.Cells.Item($r, 2).Value = $schemes17to19[$i]
    for ($c = 3; $c -le 16; $c++) {
        $ws.Cells.Item($r, $c).Value = 1
    }
}
